$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header Z1: NewCol -> SecondaryCause
$ws.Range("Z1").Value = "SecondaryCause"

# Updated SecondaryCause (formerly NewCol) categorical values per row,
# reflecting the new Air_/Env_/Per_ prefixed taxonomy.
$secondaryCause = @{}
$secondaryCause[2] = "Air_Maintenance Issue"
$secondaryCause[3] = "Env_Runway"
$secondaryCause[4] = "Env_Lighting/Visibility"
$secondaryCause[5] = "Env_Object/Animal"
$secondaryCause[6] = "Env_Runway"
$secondaryCause[7] = "Env_Object/Animal"
$secondaryCause[8] = "Air_Failure"
$secondaryCause[9] = "Aircraft Oper/Perf"
$secondaryCause[10] = "Aircraft Oper/Perf"
$secondaryCause[11] = "Pilot"
$secondaryCause[12] = "Env_Wind/Weather"
$secondaryCause[13] = "Aircraft Oper/Perf"
$secondaryCause[14] = "Pilot"
$secondaryCause[15] = "Pilot"
$secondaryCause[16] = "Pilot"
$secondaryCause[17] = "Air_Failure"
$secondaryCause[18] = "Aircraft Oper/Perf"
$secondaryCause[19] = "Aircraft Oper/Perf"
$secondaryCause[20] = "Aircraft Oper/Perf"
$secondaryCause[21] = "Crew"
$secondaryCause[22] = "Aircraft Oper/Perf"
$secondaryCause[23] = "Pilot"
$secondaryCause[24] = "Pilot"
$secondaryCause[25] = "Aircraft Oper/Perf"
$secondaryCause[26] = "Pilot"
$secondaryCause[27] = "Pilot"
$secondaryCause[28] = "Aircraft Oper/Perf"
$secondaryCause[29] = "Air_Failure"
$secondaryCause[30] = "Env_Other"
$secondaryCause[31] = "Aircraft Oper/Perf"
$secondaryCause[32] = "Pilot"
$secondaryCause[33] = "Air_Failure"
$secondaryCause[34] = "Aircraft Oper/Perf"
$secondaryCause[35] = "Aircraft Oper/Perf"
$secondaryCause[36] = "Air_Damaged Part"
$secondaryCause[37] = "Pilot"
$secondaryCause[38] = "Not Determined"
$secondaryCause[39] = "Aircraft Oper/Perf"
$secondaryCause[40] = "Env_Turbulence"
$secondaryCause[41] = "Aircraft Oper/Perf"
$secondaryCause[42] = "Aircraft Oper/Perf"
$secondaryCause[43] = "Air_Damaged Part"
$secondaryCause[44] = "Aircraft Oper/Perf"
$secondaryCause[45] = "Air_Failure"
$secondaryCause[46] = "Env_Wind/Weather"
$secondaryCause[47] = "Not Determined"
$secondaryCause[48] = "Not Determined"
$secondaryCause[49] = "Aircraft Oper/Perf"
$secondaryCause[50] = "Env_Wind/Weather"
$secondaryCause[51] = "Air_Failure"
$secondaryCause[52] = "Env_Runway"
$secondaryCause[53] = "Env_Runway"
$secondaryCause[54] = "Pilot"
$secondaryCause[55] = "Air_Damaged Part"
$secondaryCause[56] = "Aircraft Oper/Perf"
$secondaryCause[57] = "Pilot"
$secondaryCause[58] = "Pilot"
$secondaryCause[59] = "Env_Other"
$secondaryCause[60] = "Air_Damaged Part"
$secondaryCause[61] = "Pilot"
$secondaryCause[62] = "Aircraft Oper/Perf"
$secondaryCause[63] = "Pilot"
$secondaryCause[64] = "Air_Other"
$secondaryCause[65] = "Pilot"
$secondaryCause[66] = "Pilot"
$secondaryCause[67] = "Air_Failure"
$secondaryCause[68] = "Pilot"
$secondaryCause[69] = "Pilot"
$secondaryCause[70] = "Pilot"
$secondaryCause[71] = "Pilot"
$secondaryCause[72] = "Pilot"
$secondaryCause[73] = "Env_Other"
$secondaryCause[74] = "Pilot"
$secondaryCause[75] = "Pilot"
$secondaryCause[76] = "Pilot"
$secondaryCause[77] = "Pilot"
$secondaryCause[78] = "Pilot"
$secondaryCause[79] = "Pilot"
$secondaryCause[80] = "Air_Failure"
$secondaryCause[81] = "Air_Failure"
$secondaryCause[82] = "Air_Damaged Part"
$secondaryCause[83] = "Pilot"
$secondaryCause[84] = "Passenger"
$secondaryCause[85] = "Pilot"
$secondaryCause[86] = "Pilot"
$secondaryCause[87] = "Pilot"
$secondaryCause[88] = "Air_Damaged Part"
$secondaryCause[89] = "Pilot"
$secondaryCause[90] = "Aircraft Oper/Perf"
$secondaryCause[91] = "Air_Failure"
$secondaryCause[92] = "Crew"
$secondaryCause[93] = "Air_Failure"
$secondaryCause[94] = "Pilot"
$secondaryCause[95] = "Air_Failure"
$secondaryCause[96] = "Aircraft Oper/Perf"
$secondaryCause[97] = "Pilot"
$secondaryCause[98] = "Pilot"
$secondaryCause[99] = "Org Issues"
$secondaryCause[100] = "Pilot"
$secondaryCause[101] = "Air_Failure"
$secondaryCause[102] = "Aircraft Oper/Perf"
$secondaryCause[103] = "Air_Failure"
$secondaryCause[104] = "Pilot"
$secondaryCause[105] = "Pilot"
$secondaryCause[106] = "Aircraft Oper/Perf"
$secondaryCause[107] = "Aircraft Oper/Perf"
$secondaryCause[108] = "Crew"
$secondaryCause[109] = "Pilot"
$secondaryCause[110] = "Pilot"
$secondaryCause[111] = "Pilot"
$secondaryCause[112] = "Pilot"
$secondaryCause[113] = "Air_Failure"
$secondaryCause[114] = "Pilot"
$secondaryCause[115] = "Pilot"
$secondaryCause[116] = "Pilot"
$secondaryCause[117] = "Air_Failure"
$secondaryCause[118] = "Air_Maintenance Issue"
$secondaryCause[119] = "Pilot"
$secondaryCause[120] = "Not Determined"
$secondaryCause[121] = "Env_Wind/Weather"
$secondaryCause[122] = "Pilot"
$secondaryCause[123] = "Pilot"
$secondaryCause[124] = "Air_Failure"
$secondaryCause[125] = "Pilot"
$secondaryCause[126] = "Pilot"
$secondaryCause[127] = "Pilot"
$secondaryCause[128] = "Pilot"
$secondaryCause[129] = "Not Determined"
$secondaryCause[130] = "Aircraft Oper/Perf"
$secondaryCause[131] = "Pilot"
$secondaryCause[132] = "Pilot"
$secondaryCause[133] = "Pilot"
$secondaryCause[134] = "Pilot"
$secondaryCause[135] = "Pilot"
$secondaryCause[136] = "Pilot"
$secondaryCause[137] = "Pilot"
$secondaryCause[138] = "Pilot"
$secondaryCause[139] = "Air_Failure"
$secondaryCause[140] = "Pilot"
$secondaryCause[141] = "Crew"
$secondaryCause[142] = "Air_Other"
$secondaryCause[143] = "Env_Runway"
$secondaryCause[144] = "Pilot"
$secondaryCause[145] = "Pilot"
$secondaryCause[146] = "Pilot"
$secondaryCause[147] = "Env_Object/Animal"
$secondaryCause[148] = "Not Determined"
$secondaryCause[149] = "Pilot"
$secondaryCause[150] = "Pilot"
$secondaryCause[151] = "Pilot"
$secondaryCause[152] = "Pilot"
$secondaryCause[153] = "Aircraft Oper/Perf"
$secondaryCause[154] = "Aircraft Oper/Perf"
$secondaryCause[155] = "Pilot"
$secondaryCause[156] = "Aircraft Oper/Perf"
$secondaryCause[157] = "Pilot"
$secondaryCause[158] = "Aircraft Oper/Perf"
$secondaryCause[159] = "Pilot"
$secondaryCause[160] = "Pilot"
$secondaryCause[161] = "Aircraft Oper/Perf"
$secondaryCause[162] = "Aircraft Oper/Perf"
$secondaryCause[163] = "Aircraft Oper/Perf"
$secondaryCause[164] = "Env_Lighting/Visibility"
$secondaryCause[165] = "Aircraft Oper/Perf"
$secondaryCause[166] = "Pilot"
$secondaryCause[167] = "Pilot"
$secondaryCause[168] = "Pilot"
$secondaryCause[169] = "Pilot"
$secondaryCause[170] = "Not Determined"
$secondaryCause[171] = "Air_Damaged Part"
$secondaryCause[172] = "Pilot"
$secondaryCause[173] = "Aircraft Oper/Perf"
$secondaryCause[174] = "Air_Failure"
$secondaryCause[175] = "Pilot"
$secondaryCause[176] = "Pilot"
$secondaryCause[177] = "Pilot"
$secondaryCause[178] = "Env_Object/Animal"
$secondaryCause[179] = "Air_Damaged Part"
$secondaryCause[180] = "Pilot"
$secondaryCause[181] = "Pilot"
$secondaryCause[182] = "Pilot"
$secondaryCause[183] = "Pilot"
$secondaryCause[184] = "Pilot"
$secondaryCause[185] = "Air_Failure"
$secondaryCause[186] = "Pilot"
$secondaryCause[187] = "Pilot"
$secondaryCause[188] = "Per_Other"
$secondaryCause[189] = "Pilot"
$secondaryCause[190] = "Aircraft Oper/Perf"
$secondaryCause[191] = "Pilot"
$secondaryCause[192] = "Pilot"
$secondaryCause[193] = "Pilot"
$secondaryCause[194] = "Pilot"
$secondaryCause[195] = "Pilot"
$secondaryCause[196] = "Air_Other"
$secondaryCause[197] = "Aircraft Oper/Perf"
$secondaryCause[198] = "Aircraft Oper/Perf"
$secondaryCause[199] = "Aircraft Oper/Perf"
$secondaryCause[200] = "Pilot"
$secondaryCause[201] = "Aircraft Oper/Perf"
$secondaryCause[202] = "Pilot"
$secondaryCause[203] = "Aircraft Oper/Perf"
$secondaryCause[204] = "Pilot"
$secondaryCause[205] = "Aircraft Oper/Perf"
$secondaryCause[206] = "Aircraft Oper/Perf"
$secondaryCause[207] = "Airport/Airline Personnel"
$secondaryCause[208] = "Air_Failure"
$secondaryCause[209] = "Env_Lighting/Visibility"
$secondaryCause[210] = "Env_Lighting/Visibility"
$secondaryCause[211] = "Env_Turbulence"
$secondaryCause[212] = "Aircraft Oper/Perf"
$secondaryCause[213] = "Env_Object/Animal"
$secondaryCause[214] = "Pilot"
$secondaryCause[215] = "Crew"
$secondaryCause[216] = "Pilot"
$secondaryCause[217] = "Crew"
$secondaryCause[218] = "Air_Maintenance Issue"
$secondaryCause[219] = "Crew"
$secondaryCause[220] = "Env_Lighting/Visibility"
$secondaryCause[221] = "Aircraft Oper/Perf"
$secondaryCause[222] = "Pilot"
$secondaryCause[223] = "Crew"
$secondaryCause[224] = "Air_Other"
$secondaryCause[225] = "Aircraft Oper/Perf"
$secondaryCause[226] = "Env_Object/Animal"
$secondaryCause[227] = "Ground Crew"
$secondaryCause[228] = "Env_Object/Animal"
$secondaryCause[229] = "Pilot"
$secondaryCause[230] = "Air_Failure"
$secondaryCause[231] = "Ground Crew"
$secondaryCause[232] = "Aircraft Oper/Perf"
$secondaryCause[233] = "Env_Turbulence"
$secondaryCause[234] = "Env_Turbulence"
$secondaryCause[235] = "Ground Crew"
$secondaryCause[236] = "Crew"
$secondaryCause[237] = "Pilot"
$secondaryCause[238] = "Pilot"
$secondaryCause[239] = "Env_Turbulence"
$secondaryCause[240] = "Air_Damaged Part"
$secondaryCause[241] = "Crew"
$secondaryCause[242] = "Aircraft Oper/Perf"
$secondaryCause[243] = "Env_Turbulence"
$secondaryCause[244] = "Ground Crew"
$secondaryCause[245] = "Aircraft Oper/Perf"
$secondaryCause[246] = "Air_Damaged Part"
$secondaryCause[247] = "Air_Failure"
$secondaryCause[248] = "Crew"
$secondaryCause[249] = "Env_Object/Animal"
$secondaryCause[250] = "Crew"
$secondaryCause[251] = "Passenger"
$secondaryCause[252] = "Aircraft Oper/Perf"
$secondaryCause[253] = "Ground Crew"
$secondaryCause[254] = "Org Issues"
$secondaryCause[255] = "Aircraft Oper/Perf"
$secondaryCause[256] = "Env_Turbulence"
$secondaryCause[257] = "Pilot"
$secondaryCause[258] = "Env_Object/Animal"
$secondaryCause[259] = "Crew"
$secondaryCause[260] = "Pilot"
$secondaryCause[261] = "Pilot"
$secondaryCause[262] = "Air_Other"
$secondaryCause[263] = "Pilot"
$secondaryCause[264] = "Air_Damaged Part"
$secondaryCause[265] = "Air_Failure"
$secondaryCause[266] = "Pilot"
$secondaryCause[267] = "Aircraft Oper/Perf"
$secondaryCause[268] = "Env_Turbulence"
$secondaryCause[269] = "Pilot"
$secondaryCause[270] = "Pilot"
$secondaryCause[271] = "Air_Failure"
$secondaryCause[272] = "Env_Turbulence"
$secondaryCause[273] = "Air_Failure"
$secondaryCause[274] = "Pilot"
$secondaryCause[275] = "Ground Crew"
$secondaryCause[276] = "Ground Crew"
$secondaryCause[277] = "Env_Turbulence"
$secondaryCause[278] = "Pilot"
$secondaryCause[279] = "Air_Failure"
$secondaryCause[280] = "Crew"
$secondaryCause[281] = "Pilot"
$secondaryCause[282] = "Env_Turbulence"
$secondaryCause[283] = "Not Determined"
$secondaryCause[284] = "Env_Turbulence"
$secondaryCause[285] = "Env_Wind/Weather"
$secondaryCause[286] = "Env_Turbulence"
$secondaryCause[287] = "Env_Turbulence"
$secondaryCause[288] = "Env_Object/Animal"
$secondaryCause[289] = "Pilot"
$secondaryCause[290] = "Ground Crew"
$secondaryCause[291] = "Env_Turbulence"
$secondaryCause[292] = "Crew"
$secondaryCause[293] = "Ground Crew"
$secondaryCause[294] = "Env_Turbulence"
$secondaryCause[295] = "Env_Turbulence"
$secondaryCause[296] = "Env_Turbulence"
$secondaryCause[297] = "Air_Failure"
$secondaryCause[298] = "Aircraft Oper/Perf"
$secondaryCause[299] = "Aircraft Oper/Perf"
$secondaryCause[300] = "Air_Other"
$secondaryCause[301] = "Env_Other"
$secondaryCause[302] = "Env_Turbulence"
$secondaryCause[303] = "Air_Failure"
$secondaryCause[304] = "Aircraft Oper/Perf"
$secondaryCause[305] = "Air_Damaged Part"
$secondaryCause[306] = "Crew"
$secondaryCause[307] = "Env_Turbulence"
$secondaryCause[308] = "Pilot"
$secondaryCause[309] = "Air_Maintenance Issue"
$secondaryCause[310] = "Per_Other"
$secondaryCause[311] = "Env_Turbulence"
$secondaryCause[312] = "Ground Crew"
$secondaryCause[313] = "Pilot"
$secondaryCause[314] = "Env_Turbulence"
$secondaryCause[315] = "Env_Runway"
$secondaryCause[316] = "Env_Turbulence"
$secondaryCause[317] = "Env_Wind/Weather"
$secondaryCause[318] = "Crew"
$secondaryCause[319] = "Crew"
$secondaryCause[320] = "Env_Object/Animal"
$secondaryCause[321] = "Env_Turbulence"
$secondaryCause[322] = "Crew"
$secondaryCause[323] = "Airport/Airline Personnel"
$secondaryCause[324] = "Pilot"
$secondaryCause[325] = "Env_Object/Animal"
$secondaryCause[326] = "Env_Wind/Weather"
$secondaryCause[327] = "Aircraft Oper/Perf"
$secondaryCause[328] = "Crew"
$secondaryCause[329] = "Air_Failure"
$secondaryCause[330] = "Env_Equipment"
$secondaryCause[331] = "Ground Crew"
$secondaryCause[332] = "Ground Crew"
$secondaryCause[333] = "Pilot"
$secondaryCause[334] = "Env_Wind/Weather"
$secondaryCause[335] = "Ground Crew"
$secondaryCause[336] = "Air_Failure"
$secondaryCause[337] = "Crew"
$secondaryCause[338] = "Air_Maintenance Issue"
$secondaryCause[339] = "Ground Crew"
$secondaryCause[340] = "Crew"
$secondaryCause[341] = "Env_Turbulence"
$secondaryCause[342] = "Ground Crew"
$secondaryCause[343] = "Env_Other"
$secondaryCause[344] = "Env_Turbulence"
$secondaryCause[345] = "Pilot"
$secondaryCause[346] = "Env_Equipment"
$secondaryCause[347] = "Env_Turbulence"
$secondaryCause[348] = "Air_Damaged Part"
$secondaryCause[349] = "Pilot"
$secondaryCause[350] = "Pilot"
$secondaryCause[351] = "Env_Turbulence"
$secondaryCause[352] = "Crew"
$secondaryCause[353] = "Env_Turbulence"
$secondaryCause[354] = "Pilot"
$secondaryCause[355] = "Air_Other"
$secondaryCause[356] = "Pilot"
$secondaryCause[357] = "Crew"
$secondaryCause[358] = "Env_Turbulence"
$secondaryCause[359] = "Env_Turbulence"
$secondaryCause[360] = "Airport/Airline Personnel"
$secondaryCause[361] = "Env_Turbulence"
$secondaryCause[362] = "Crew"
$secondaryCause[363] = "Crew"
$secondaryCause[364] = "Pilot"
$secondaryCause[365] = "Pilot"
$secondaryCause[366] = "Pilot"
$secondaryCause[367] = "Air_Failure"
$secondaryCause[368] = "Pilot"
$secondaryCause[369] = "Ground Crew"
$secondaryCause[370] = "Env_Turbulence"
$secondaryCause[371] = "Pilot"
$secondaryCause[372] = "Air_Failure"
$secondaryCause[373] = "Env_Runway"
$secondaryCause[374] = "Env_Turbulence"
$secondaryCause[375] = "Pilot"
$secondaryCause[376] = "Aircraft Oper/Perf"
$secondaryCause[377] = "Air_Failure"
$secondaryCause[378] = "Crew"
$secondaryCause[379] = "Env_Turbulence"
$secondaryCause[380] = "Env_Turbulence"
$secondaryCause[381] = "Env_Turbulence"
$secondaryCause[382] = "Air_Damaged Part"
$secondaryCause[383] = "Air_Damaged Part"
$secondaryCause[384] = "Air_Maintenance Issue"
$secondaryCause[385] = "Env_Turbulence"
$secondaryCause[386] = "Env_Turbulence"
$secondaryCause[387] = "Env_Turbulence"
$secondaryCause[388] = "Air_Failure"
$secondaryCause[389] = "Pilot"
$secondaryCause[390] = "Pilot"
$secondaryCause[391] = "Not Determined"
$secondaryCause[392] = "Crew"
$secondaryCause[393] = "Air_Maintenance Issue"
$secondaryCause[394] = "Env_Turbulence"
$secondaryCause[395] = "Air_Damaged Part"
$secondaryCause[396] = "Env_Turbulence"
$secondaryCause[397] = "Pilot"
$secondaryCause[398] = "Crew"
$secondaryCause[399] = "Env_Turbulence"
$secondaryCause[400] = "Air_Failure"
$secondaryCause[401] = "Env_Turbulence"
$secondaryCause[402] = "Airport/Airline Personnel"
$secondaryCause[403] = "Env_Equipment"
$secondaryCause[404] = "Env_Wind/Weather"
$secondaryCause[405] = "Env_Turbulence"
$secondaryCause[406] = "Airport/Airline Personnel"
$secondaryCause[407] = "Pilot"
$secondaryCause[408] = "Airport/Airline Personnel"
$secondaryCause[409] = "Pilot"
$secondaryCause[410] = "Env_Turbulence"
$secondaryCause[411] = "Air_Other"
$secondaryCause[412] = "Air_Other"
$secondaryCause[413] = "Not Determined"
$secondaryCause[414] = "Pilot"
$secondaryCause[415] = "Env_Wind/Weather"
$secondaryCause[416] = "Airport/Airline Personnel"
$secondaryCause[417] = "Env_Equipment"
$secondaryCause[418] = "Aircraft Oper/Perf"
$secondaryCause[419] = "Air_Damaged Part"
$secondaryCause[420] = "Env_Turbulence"
$secondaryCause[421] = "Pilot"
$secondaryCause[422] = "Pilot"
$secondaryCause[423] = "Env_Turbulence"
$secondaryCause[424] = "Air_Failure"
$secondaryCause[425] = "Env_Turbulence"
$secondaryCause[426] = "Env_Turbulence"
$secondaryCause[427] = "Pilot"
$secondaryCause[428] = "Air_Maintenance Issue"
$secondaryCause[429] = "Air_Damaged Part"
$secondaryCause[430] = "Air_Failure"
$secondaryCause[431] = "Airport/Airline Personnel"
$secondaryCause[432] = "Air_Failure"
$secondaryCause[433] = "Air_Maintenance Issue"
$secondaryCause[434] = "Air_Failure"
$secondaryCause[435] = "Aircraft Oper/Perf"
$secondaryCause[436] = "Pilot"
$secondaryCause[437] = "Aircraft Oper/Perf"
$secondaryCause[438] = "Aircraft Oper/Perf"
$secondaryCause[439] = "Air_Other"
$secondaryCause[440] = "Aircraft Oper/Perf"
$secondaryCause[441] = "Pilot"
$secondaryCause[442] = "Ground Crew"
$secondaryCause[443] = "Passenger"
$secondaryCause[444] = "Crew"
$secondaryCause[445] = "Env_Equipment"
$secondaryCause[446] = "Crew"
$secondaryCause[447] = "Passenger"
$secondaryCause[448] = "Env_Object/Animal"
$secondaryCause[449] = "Crew"
$secondaryCause[450] = "Aircraft Oper/Perf"
$secondaryCause[451] = "Pilot"
$secondaryCause[452] = "Aircraft Oper/Perf"
$secondaryCause[453] = "Pilot"
$secondaryCause[454] = "Crew"
$secondaryCause[455] = "Env_Turbulence"
$secondaryCause[456] = "Per_Other"
$secondaryCause[457] = "Pilot"
$secondaryCause[458] = "Env_Turbulence"
$secondaryCause[459] = "Pilot"
$secondaryCause[460] = "Pilot"
$secondaryCause[461] = "Pilot"
$secondaryCause[462] = "Env_Turbulence"
$secondaryCause[463] = "Pilot"
$secondaryCause[464] = "Ground Crew"
$secondaryCause[465] = "Aircraft Oper/Perf"
$secondaryCause[466] = "Ground Crew"
$secondaryCause[467] = "Pilot"
$secondaryCause[468] = "Env_Turbulence"
$secondaryCause[469] = "Pilot"
$secondaryCause[470] = "Env_Object/Animal"
$secondaryCause[471] = "Aircraft Oper/Perf"
$secondaryCause[472] = "Pilot"
$secondaryCause[473] = "Not Determined"
$secondaryCause[474] = "Env_Wind/Weather"
$secondaryCause[475] = "Pilot"
$secondaryCause[476] = "Env_Turbulence"
$secondaryCause[477] = "Air_Failure"
$secondaryCause[478] = "Env_Turbulence"
$secondaryCause[479] = "Aircraft Oper/Perf"
$secondaryCause[480] = "Aircraft Oper/Perf"
$secondaryCause[481] = "Env_Turbulence"
$secondaryCause[482] = "Air_Damaged Part"
$secondaryCause[483] = "Env_Turbulence"
$secondaryCause[484] = "Per_Other"
$secondaryCause[485] = "Aircraft Oper/Perf"
$secondaryCause[486] = "Pilot"
$secondaryCause[487] = "Aircraft Oper/Perf"
$secondaryCause[488] = "Env_Wind/Weather"
$secondaryCause[489] = "Crew"
$secondaryCause[490] = "Crew"
$secondaryCause[491] = "Pilot"
$secondaryCause[492] = "Crew"
$secondaryCause[493] = "Env_Turbulence"
$secondaryCause[494] = "Pilot"
$secondaryCause[495] = "Air_Failure"
$secondaryCause[496] = "Ground Crew"
$secondaryCause[497] = "Env_Turbulence"
$secondaryCause[498] = "Pilot"
$secondaryCause[499] = "Air_Failure"
$secondaryCause[500] = "Per_Other"
$secondaryCause[501] = "Env_Turbulence"
$secondaryCause[502] = "Env_Turbulence"
$secondaryCause[503] = "Env_Runway"
$secondaryCause[504] = "Env_Turbulence"
$secondaryCause[505] = "Env_Turbulence"
$secondaryCause[506] = "Aircraft Oper/Perf"
$secondaryCause[507] = "Env_Turbulence"
$secondaryCause[508] = "Aircraft Oper/Perf"
$secondaryCause[509] = "Air_Other"
$secondaryCause[510] = "Env_Turbulence"
$secondaryCause[511] = "Passenger"
$secondaryCause[512] = "Env_Turbulence"
$secondaryCause[513] = "Passenger"
$secondaryCause[514] = "Pilot"
$secondaryCause[515] = "Air_Failure"
$secondaryCause[516] = "Env_Turbulence"
$secondaryCause[517] = "Env_Turbulence"
$secondaryCause[518] = "Pilot"
$secondaryCause[519] = "Not Determined"
$secondaryCause[520] = "Aircraft Oper/Perf"
$secondaryCause[521] = "Env_Runway"
$secondaryCause[522] = "Env_Object/Animal"
$secondaryCause[523] = "Air_Other"
$secondaryCause[524] = "Aircraft Oper/Perf"
$secondaryCause[525] = "Ground Crew"
$secondaryCause[526] = "Crew"
$secondaryCause[527] = "Not Determined"
$secondaryCause[528] = "Pilot"
$secondaryCause[529] = "Pilot"
$secondaryCause[530] = "Pilot"
$secondaryCause[531] = "Env_Object/Animal"
$secondaryCause[532] = "Air_Other"
$secondaryCause[533] = "Env_Lighting/Visibility"
$secondaryCause[534] = "Air_Damaged Part"
$secondaryCause[535] = "Env_Turbulence"
$secondaryCause[536] = "Pilot"
$secondaryCause[537] = "Pilot"
$secondaryCause[538] = "Pilot"
$secondaryCause[539] = "Crew"
$secondaryCause[540] = "Pilot"
$secondaryCause[541] = "Aircraft Oper/Perf"
$secondaryCause[542] = "Aircraft Oper/Perf"
$secondaryCause[543] = "Pilot"
$secondaryCause[544] = "Env_Turbulence"
$secondaryCause[545] = "Pilot"
$secondaryCause[546] = "Air_Failure"
$secondaryCause[547] = "Pilot"
$secondaryCause[548] = "Pilot"
$secondaryCause[549] = "Passenger"
$secondaryCause[550] = "Env_Turbulence"
$secondaryCause[551] = "Env_Object/Animal"
$secondaryCause[552] = "Crew"
$secondaryCause[553] = "Env_Lighting/Visibility"
$secondaryCause[554] = "Pilot"
$secondaryCause[555] = "Env_Object/Animal"
$secondaryCause[556] = "Env_Object/Animal"
$secondaryCause[557] = "Env_Turbulence"
$secondaryCause[558] = "Env_Runway"
$secondaryCause[559] = "Pilot"
$secondaryCause[560] = "Passenger"
$secondaryCause[561] = "Aircraft Oper/Perf"
$secondaryCause[562] = "Env_Turbulence"
$secondaryCause[563] = "Env_Turbulence"
$secondaryCause[564] = "Aircraft Oper/Perf"
$secondaryCause[565] = "Ground Crew"
$secondaryCause[566] = "Env_Turbulence"
$secondaryCause[567] = "Pilot"
$secondaryCause[568] = "Pilot"
$secondaryCause[569] = "Air_Maintenance Issue"
$secondaryCause[570] = "Env_Turbulence"
$secondaryCause[571] = "Env_Turbulence"
$secondaryCause[572] = "Env_Object/Animal"
$secondaryCause[573] = "Env_Turbulence"
$secondaryCause[574] = "Crew"
$secondaryCause[575] = "Crew"
$secondaryCause[576] = "Env_Object/Animal"
$secondaryCause[577] = "Air_Damaged Part"
$secondaryCause[578] = "Not Determined"
$secondaryCause[579] = "Env_Runway"
$secondaryCause[580] = "Crew"
$secondaryCause[581] = "Aircraft Oper/Perf"
$secondaryCause[582] = "Aircraft Oper/Perf"
$secondaryCause[583] = "Env_Object/Animal"
$secondaryCause[584] = "Crew"
$secondaryCause[585] = "Air_Damaged Part"
$secondaryCause[586] = "Aircraft Oper/Perf"
$secondaryCause[587] = "Pilot"
$secondaryCause[588] = "Pilot"
$secondaryCause[589] = "Air_Failure"
$secondaryCause[590] = "Pilot"
$secondaryCause[591] = "Ground Crew"
$secondaryCause[592] = "Air_Damaged Part"
$secondaryCause[593] = "Pilot"
$secondaryCause[594] = "Crew"
$secondaryCause[595] = "Env_Turbulence"
$secondaryCause[596] = "Aircraft Oper/Perf"
$secondaryCause[597] = "Crew"
$secondaryCause[598] = "Org Issues"
$secondaryCause[599] = "Air_Other"
$secondaryCause[600] = "Aircraft Oper/Perf"
$secondaryCause[601] = "Passenger"
$secondaryCause[602] = "Env_Turbulence"
$secondaryCause[603] = "Env_Object/Animal"
$secondaryCause[604] = "Pilot"
$secondaryCause[605] = "Air_Maintenance Issue"
$secondaryCause[606] = "Env_Turbulence"
$secondaryCause[607] = "Passenger"
$secondaryCause[608] = "Air_Damaged Part"
$secondaryCause[609] = "Air_Failure"
$secondaryCause[610] = "Env_Turbulence"
$secondaryCause[611] = "Aircraft Oper/Perf"
$secondaryCause[612] = "Env_Object/Animal"
$secondaryCause[613] = "Aircraft Oper/Perf"
$secondaryCause[614] = "Pilot"
$secondaryCause[615] = "Env_Turbulence"
$secondaryCause[616] = "Aircraft Oper/Perf"
$secondaryCause[617] = "Env_Object/Animal"
$secondaryCause[618] = "Per_Other"
$secondaryCause[619] = "Env_Turbulence"
$secondaryCause[620] = "Air_Other"
$secondaryCause[621] = "Air_Failure"
$secondaryCause[622] = "Crew"
$secondaryCause[623] = "Ground Crew"
$secondaryCause[624] = "Ground Crew"
$secondaryCause[625] = "Pilot"
$secondaryCause[626] = "Crew"
$secondaryCause[627] = "Aircraft Oper/Perf"
$secondaryCause[628] = "Aircraft Oper/Perf"
$secondaryCause[629] = "Aircraft Oper/Perf"
$secondaryCause[630] = "Aircraft Oper/Perf"
$secondaryCause[631] = "Env_Wind/Weather"
$secondaryCause[632] = "Air_Failure"
$secondaryCause[633] = "Not Determined"
$secondaryCause[634] = "Air_Failure"
$secondaryCause[635] = "Pilot"
$secondaryCause[636] = "Pilot"
$secondaryCause[637] = "Env_Turbulence"
$secondaryCause[638] = "Pilot"
$secondaryCause[639] = "Air_Damaged Part"
$secondaryCause[640] = "Env_Turbulence"
$secondaryCause[641] = "Air_Failure"
$secondaryCause[642] = "Aircraft Oper/Perf"
$secondaryCause[643] = "Env_Wind/Weather"
$secondaryCause[644] = "Air_Damaged Part"
$secondaryCause[645] = "Env_Turbulence"
$secondaryCause[646] = "Airport/Airline Personnel"
$secondaryCause[647] = "Env_Equipment"
$secondaryCause[648] = "Aircraft Oper/Perf"
$secondaryCause[649] = "Env_Turbulence"
$secondaryCause[650] = "Pilot"
$secondaryCause[651] = "Crew"
$secondaryCause[652] = "Org Issues"
$secondaryCause[653] = "Airport/Airline Personnel"
$secondaryCause[654] = "Crew"
$secondaryCause[655] = "Env_Turbulence"
$secondaryCause[656] = "Pilot"
$secondaryCause[657] = "Env_Turbulence"
$secondaryCause[658] = "Env_Turbulence"
$secondaryCause[659] = "Not Determined"
$secondaryCause[660] = "Pilot"
$secondaryCause[661] = "Aircraft Oper/Perf"
$secondaryCause[662] = "Crew"
$secondaryCause[663] = "Env_Wind/Weather"
$secondaryCause[664] = "Aircraft Oper/Perf"
$secondaryCause[665] = "Airport/Airline Personnel"
$secondaryCause[666] = "Pilot"
$secondaryCause[667] = "Crew"
$secondaryCause[668] = "Env_Turbulence"
$secondaryCause[669] = "Env_Turbulence"
$secondaryCause[670] = "Env_Turbulence"
$secondaryCause[671] = "Env_Turbulence"
$secondaryCause[672] = "Not Determined"
$secondaryCause[673] = "Aircraft Oper/Perf"
$secondaryCause[674] = "Env_Object/Animal"
$secondaryCause[675] = "Env_Turbulence"
$secondaryCause[676] = "Env_Object/Animal"
$secondaryCause[677] = "Crew"
$secondaryCause[678] = "Env_Object/Animal"
$secondaryCause[679] = "Aircraft Oper/Perf"
$secondaryCause[680] = "Ground Crew"
$secondaryCause[681] = "Ground Crew"
$secondaryCause[682] = "Env_Turbulence"
$secondaryCause[683] = "Env_Object/Animal"
$secondaryCause[684] = "Env_Object/Animal"
$secondaryCause[685] = "Pilot"
$secondaryCause[686] = "Pilot"
$secondaryCause[687] = "Env_Object/Animal"
$secondaryCause[688] = "Env_Turbulence"
$secondaryCause[689] = "Env_Turbulence"
$secondaryCause[690] = "Aircraft Oper/Perf"
$secondaryCause[691] = "Env_Object/Animal"
$secondaryCause[692] = "Air_Damaged Part"
$secondaryCause[693] = "Pilot"
$secondaryCause[694] = "Ground Crew"
$secondaryCause[695] = "Env_Turbulence"
$secondaryCause[696] = "Pilot"
$secondaryCause[697] = "Aircraft Oper/Perf"
$secondaryCause[698] = "Env_Wind/Weather"
$secondaryCause[699] = "Env_Object/Animal"
$secondaryCause[700] = "Aircraft Oper/Perf"
$secondaryCause[701] = "Aircraft Oper/Perf"
$secondaryCause[702] = "Crew"
$secondaryCause[703] = "Env_Turbulence"
$secondaryCause[704] = "Env_Lighting/Visibility"
$secondaryCause[705] = "Pilot"
$secondaryCause[706] = "Pilot"
$secondaryCause[707] = "Pilot"
$secondaryCause[708] = "Air_Failure"
$secondaryCause[709] = "Aircraft Oper/Perf"
$secondaryCause[710] = "Air_Failure"
$secondaryCause[711] = "Aircraft Oper/Perf"
$secondaryCause[712] = "Env_Runway"
$secondaryCause[713] = "Env_Object/Animal"
$secondaryCause[714] = "Pilot"
$secondaryCause[715] = "Env_Turbulence"
$secondaryCause[716] = "Env_Turbulence"
$secondaryCause[717] = "Env_Turbulence"
$secondaryCause[718] = "Env_Object/Animal"
$secondaryCause[719] = "Env_Object/Animal"
$secondaryCause[720] = "Aircraft Oper/Perf"
$secondaryCause[721] = "Ground Crew"
$secondaryCause[722] = "Pilot"
$secondaryCause[723] = "Air_Failure"
$secondaryCause[724] = "Aircraft Oper/Perf"
$secondaryCause[725] = "Env_Turbulence"
$secondaryCause[726] = "Pilot"
$secondaryCause[727] = "Pilot"
$secondaryCause[728] = "Air_Failure"
$secondaryCause[729] = "Crew"
$secondaryCause[730] = "Pilot"
$secondaryCause[731] = "Env_Other"
$secondaryCause[732] = "Env_Turbulence"
$secondaryCause[733] = "Crew"
$secondaryCause[734] = "Env_Turbulence"
$secondaryCause[735] = "Env_Object/Animal"
$secondaryCause[736] = "Aircraft Oper/Perf"
$secondaryCause[737] = "Crew"

foreach ($row in $secondaryCause.Keys) {
    $ws.Cells.Item($row, 26).Value = $secondaryCause[$row]
}
